{"js": "// \"fixed all list and new tag associations\"\n// This edit adjusts several custom paragraph/character styles in the\n// document's style catalog (List 1, List 6, List 7, List 7 Char, List 8,\n// List 3_change, List 4_change):\n//   - drop \"contextual spacing\" (w:contextualSpacing) from List 1, List 6,\n//     List 8, List 3_change, List 4_change (List 7 keeps it)\n//   - unlink List 6 / List 7 / List 8 from their \"Heading 4\" base style\n//     and instead give them their own explicit font (Times New Roman,\n//     12pt / 12pt complex-script) formerly inherited from Heading 4\n//   - List 7 (and its linked List 7 Char) becomes bold; List 6 / List 8\n//     stay non-bold\n\nconst styles = context.document.getStyles();\n\nconst list1 = styles.getByNameOrNullObject(\"List 1\");\nconst list6 = styles.getByNameOrNullObject(\"List 6\");\nconst list7 = styles.getByNameOrNullObject(\"List 7\");\nconst list7Char = styles.getByNameOrNullObject(\"List 7 Char\");\nconst list8 = styles.getByNameOrNullObject(\"List 8\");\nconst list3change = styles.getByNameOrNullObject(\"List 3_change\");\nconst list4change = styles.getByNameOrNullObject(\"List 4_change\");\n\nlist1.load(\"isNullObject\");\nlist6.load(\"isNullObject\");\nlist7.load(\"isNullObject\");\nlist7Char.load(\"isNullObject\");\nlist8.load(\"isNullObject\");\nlist3change.load(\"isNullObject\");\nlist4change.load(\"isNullObject\");\nawait context.sync();\n\n// List 1: remove contextual spacing.\nif (!list1.isNullObject) {\n  list1.noSpaceBetweenParagraphsOfSameStyle = false;\n}\n\n// List 6: no longer based on Heading 4; loses contextual spacing; gains\n// its own explicit (non-bold) Times New Roman 12pt run formatting.\nif (!list6.isNullObject) {\n  list6.baseStyle = \"\";\n  list6.noSpaceBetweenParagraphsOfSameStyle = false;\n  list6.font.name = \"Times New Roman\";\n  list6.font.nameBidirectional = \"Times New Roman\";\n  list6.font.size = 12;\n  list6.font.sizeBidirectional = 12;\n  list6.font.bold = false;\n}\n\n// List 7: no longer based on Heading 4 (contextual spacing is kept);\n// gains its own explicit bold Times New Roman 12pt run formatting.\nif (!list7.isNullObject) {\n  list7.baseStyle = \"\";\n  list7.font.name = \"Times New Roman\";\n  list7.font.nameBidirectional = \"Times New Roman\";\n  list7.font.size = 12;\n  list7.font.sizeBidirectional = 12;\n  list7.font.bold = true;\n}\n\n// List 7 Char: the run formatting mirrors List 7 and also becomes bold.\nif (!list7Char.isNullObject) {\n  list7Char.font.bold = true;\n}\n\n// List 8: no longer based on Heading 4; loses contextual spacing; gains\n// its own explicit (non-bold) Times New Roman 12pt run formatting.\nif (!list8.isNullObject) {\n  list8.baseStyle = \"\";\n  list8.noSpaceBetweenParagraphsOfSameStyle = false;\n  list8.font.name = \"Times New Roman\";\n  list8.font.nameBidirectional = \"Times New Roman\";\n  list8.font.size = 12;\n  list8.font.sizeBidirectional = 12;\n  list8.font.bold = false;\n}\n\n// List 3_change / List 4_change: remove contextual spacing.\nif (!list3change.isNullObject) {\n  list3change.noSpaceBetweenParagraphsOfSameStyle = false;\n}\nif (!list4change.isNullObject) {\n  list4change.noSpaceBetweenParagraphsOfSameStyle = false;\n}\n\nawait context.sync();\n", "ps1": "# \"fixed all list and new tag associations\"\n# This edit adjusts several custom paragraph/character styles in the\n# document's style catalog (List 1, List 6, List 7, List 7 Char, List 8,\n# List 3_change, List 4_change):\n#   - drop \"contextual spacing\" from List 1, List 6, List 8, List 3_change,\n#     List 4_change (List 7 keeps its contextual spacing)\n#   - unlink List 6 / List 7 / List 8 from their \"Heading 4\" base style and\n#     instead give them their own explicit font (Times New Roman, 12pt /\n#     12pt complex-script) formerly inherited from Heading 4\n#   - List 7 (and its linked List 7 Char) becomes bold; List 6 / List 8\n#     stay non-bold\n\n$d = $word.ActiveDocument\n\nfunction Get-StyleOrNull($doc, $name) {\n    try {\n        return $doc.Styles($name)\n    } catch {\n        return $null\n    }\n}\n\n# List 1: remove contextual spacing.\n$list1 = Get-StyleOrNull $d \"List 1\"\nif ($list1) {\n    $list1.NoSpaceBetweenParagraphsOfSameStyle = $false\n}\n\n# List 6: no longer based on Heading 4; loses contextual spacing; gains its\n# own explicit (non-bold) Times New Roman 12pt run formatting.\n$list6 = Get-StyleOrNull $d \"List 6\"\nif ($list6) {\n    $list6.BaseStyle = \"\"\n    $list6.NoSpaceBetweenParagraphsOfSameStyle = $false\n    $list6.Font.Name = \"Times New Roman\"\n    $list6.Font.NameBi = \"Times New Roman\"\n    $list6.Font.Size = 12\n    $list6.Font.SizeBi = 12\n    $list6.Font.Bold = $false\n}\n\n# List 7: no longer based on Heading 4 (contextual spacing is kept); gains\n# its own explicit bold Times New Roman 12pt run formatting.\n$list7 = Get-StyleOrNull $d \"List 7\"\nif ($list7) {\n    $list7.BaseStyle = \"\"\n    $list7.Font.Name = \"Times New Roman\"\n    $list7.Font.NameBi = \"Times New Roman\"\n    $list7.Font.Size = 12\n    $list7.Font.SizeBi = 12\n    $list7.Font.Bold = $true\n}\n\n# List 7 Char: the run formatting mirrors List 7 and also becomes bold.\n$list7Char = Get-StyleOrNull $d \"List 7 Char\"\nif ($list7Char) {\n    $list7Char.Font.Bold = $true\n}\n\n# List 8: no longer based on Heading 4; loses contextual spacing; gains its\n# own explicit (non-bold) Times New Roman 12pt run formatting.\n$list8 = Get-StyleOrNull $d \"List 8\"\nif ($list8) {\n    $list8.BaseStyle = \"\"\n    $list8.NoSpaceBetweenParagraphsOfSameStyle = $false\n    $list8.Font.Name = \"Times New Roman\"\n    $list8.Font.NameBi = \"Times New Roman\"\n    $list8.Font.Size = 12\n    $list8.Font.SizeBi = 12\n    $list8.Font.Bold = $false\n}\n\n# List 3_change / List 4_change: remove contextual spacing.\n$list3change = Get-StyleOrNull $d \"List 3_change\"\nif ($list3change) {\n    $list3change.NoSpaceBetweenParagraphsOfSameStyle = $false\n}\n\n$list4change = Get-StyleOrNull $d \"List 4_change\"\nif ($list4change) {\n    $list4change.NoSpaceBetweenParagraphsOfSameStyle = $false\n}\n"}
